$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 6186
$ws.Range("I3").Value = 6430
$ws.Range("I4").Value = 1481
$ws.Range("I5").Value = 595
$ws.Range("I6").Value = 7310
$ws.Range("I7").Value = 22002

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I2").Value = 64
$ws.Range("I7").Value = 258

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I2").Value = 73
$ws.Range("I7").Value = 241

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I6").Value = 206
$ws.Range("I7").Value = 696

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I2").Value = 110
$ws.Range("I6").Value = 101
$ws.Range("I7").Value = 398

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I2").Value = 198
$ws.Range("I3").Value = 312
$ws.Range("I7").Value = 848

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("I2").Value = 71
$ws.Range("I7").Value = 192

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I6").Value = 162
$ws.Range("I7").Value = 697
$ws.Range("I8").Value = 1322
$ws.Range("I10").Value = 156
$ws.Range("I19").Value = 608
$ws.Range("I20").Value = 554
$ws.Range("I22").Value = 61
$ws.Range("I23").Value = 216
$ws.Range("I25").Value = 119
$ws.Range("I26").Value = 29
$ws.Range("I29").Value = 1349
$ws.Range("I33").Value = 1000
$ws.Range("I37").Value = 696
$ws.Range("I42").Value = 770
$ws.Range("I47").Value = 158
$ws.Range("I49").Value = 148
$ws.Range("I51").Value = 258
$ws.Range("I52").Value = 472
$ws.Range("I53").Value = 234
$ws.Range("I54").Value = 452
$ws.Range("I55").Value = 246
$ws.Range("I60").Value = 121
$ws.Range("I63").Value = 75
$ws.Range("I64").Value = 182
$ws.Range("I67").Value = 848
$ws.Range("I69").Value = 48
$ws.Range("I72").Value = 86
$ws.Range("I75").Value = 72
$ws.Range("I78").Value = 297
$ws.Range("I79").Value = 622
$ws.Range("I83").Value = 480
$ws.Range("I84").Value = 192
$ws.Range("I85").Value = 995
$ws.Range("I89").Value = 258
$ws.Range("I90").Value = 278
$ws.Range("I91").Value = 233
$ws.Range("I94").Value = 228
$ws.Range("I95").Value = 338
$ws.Range("I96").Value = 241
$ws.Range("I99").Value = 398
$ws.Range("I101").Value = 22002

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I3").Value = 176
$ws.Range("I7").Value = 480

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I6").Value = 70
$ws.Range("I7").Value = 338

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I5").Value = 41
$ws.Range("I6").Value = 317
$ws.Range("I7").Value = 1000

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("I6").Value = 90
$ws.Range("I7").Value = 148

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I6").Value = 217
$ws.Range("I7").Value = 452

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 398
$ws.Range("I7").Value = 1349

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I4").Value = 22
$ws.Range("I5").Value = 16
$ws.Range("I6").Value = 184
$ws.Range("I7").Value = 608

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("I2").Value = 52
$ws.Range("I6").Value = 46

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 285
$ws.Range("I3").Value = 384
$ws.Range("I7").Value = 995

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I3").Value = 41
$ws.Range("I7").Value = 162

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I3").Value = 238
$ws.Range("I7").Value = 770

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("I3").Value = 31
$ws.Range("I7").Value = 156

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I3").Value = 74
$ws.Range("I7").Value = 297

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I2").Value = 75
$ws.Range("I3").Value = 78
$ws.Range("I7").Value = 246

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I6").Value = 63
$ws.Range("I7").Value = 216

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("I6").Value = 17
$ws.Range("I7").Value = 48

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I2").Value = 73
$ws.Range("I3").Value = 85
$ws.Range("I6").Value = 64
$ws.Range("I7").Value = 233

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I2").Value = 182
$ws.Range("I6").Value = 185
$ws.Range("I7").Value = 622

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("I3").Value = 54
$ws.Range("I7").Value = 182

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I2").Value = 151
$ws.Range("I6").Value = 196
$ws.Range("I7").Value = 554

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I6").Value = 124
$ws.Range("I7").Value = 472

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("I6").Value = 130
$ws.Range("I7").Value = 228

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("I2").Value = 45
$ws.Range("I7").Value = 119

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("I2").Value = 37
$ws.Range("I7").Value = 158

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("I2").Value = 6
$ws.Range("I7").Value = 29

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 395
$ws.Range("I3").Value = 373
$ws.Range("I6").Value = 431
$ws.Range("I7").Value = 1322

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("I4").Value = 5
$ws.Range("I7").Value = 72

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I6").Value = 94
$ws.Range("I7").Value = 278

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I6").Value = 105
$ws.Range("I7").Value = 258

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("I6").Value = 36
$ws.Range("I7").Value = 121

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I3").Value = 48
$ws.Range("I6").Value = 109
$ws.Range("I7").Value = 234

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("I2").Value = 25
$ws.Range("I7").Value = 61

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("I3").Value = 18
$ws.Range("I6").Value = 41
$ws.Range("I7").Value = 86

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I2").Value = 230
$ws.Range("I3").Value = 216
$ws.Range("I6").Value = 183
$ws.Range("I7").Value = 697
